$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "IPC PO" (predicted output) values for rows 2..51, column C.
# These are the neural network's predictions that now replace the
# placeholder zeros left over from before the weight-handling refactor.
$newC = @(29.92086429301134,29.75576833558574,30.13415867165269,30.30169986370785,30.3981413549078,30.79484862147722,30.40182221066854,30.58925283100595,30.76726524769149,30.65553389362031,30.46319840433761,30.33826832813858,30.91421946098578,31.08334631199928,31.37386341866624,31.35382061469859,31.76592095163,32.08039714711408,32.05691973869513,31.72850874014626,32.34452537865486,32.3279649760885,33.35621045423139,33.01847596218533,33.19891958886334,33.16046126621927,33.7893257046001,34.5665080189096,34.59969205759258,35.16502067223599,34.96264882611138,35.26988327467259,36.0186020077862,36.39685792328989,37.41331571183255,37.87972272349179,38.38858216854828,39.22039050917978,39.91697248426966,40.06922978261161,39.85084667886793,41.03408562619104,41.26257785707261,41.01829683879443,40.88273289387499,41.34702269434033,41.86829395669196,42.17851202908643,43.42156435389256,43.668173480221)

$startRow = 2
for ($i = 0; $i -lt $newC.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value2 = $newC[$i]
}

# Recompute the dependent columns/cells exactly as the sheet's model does:
#   D (DELTA)    = C (IPC PO) - B (IPC RO)
#   E (DELTA^2)  = D * D
#   C52 (TOTAL)  = SUM(D2:D51)
#   E52 (TOTAL)  = SUM(E2:E51)
#   E53 (MSE)    = E52 / 50
$totalD = 0.0
$totalE = 0.0
for ($row = 2; $row -le 51; $row++) {
    $b = $ws.Cells.Item($row, 2).Value2
    $c = $ws.Cells.Item($row, 3).Value2
    $d = $c - $b
    $e = $d * $d
    $ws.Cells.Item($row, 4).Value2 = $d
    $ws.Cells.Item($row, 5).Value2 = $e
    $totalD += $d
    $totalE += $e
}

$ws.Cells.Item(52, 3).Value2 = $totalD
$ws.Cells.Item(52, 5).Value2 = $totalE
$ws.Cells.Item(53, 5).Value2 = $totalE / 50

Write-Output "Updated IPC PO predictions and recomputed DELTA/DELTA^2/TOTAL/MSE"
